$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login credentials used by the test case:
#   A2 (Username) : makaia@testleaf.com   -> matschie@testleaf.com
#   B2 (Password)  : SelBootcamp$123       -> SelBootCamp@123
$ws.Range("A2").Value = "matschie@testleaf.com"
$ws.Range("B2").Value = "SelBootCamp@123"

# Refresh the mailto hyperlinks so their targets match the new values, and
# give the Password cell a hyperlink too (mirroring the Username cell).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:matschie@testleaf.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:SelBootCamp@123") | Out-Null

# Both credential cells use the Hyperlink style now.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"

# Move the active selection to A2.
$ws.Range("A2").Select()
